# Fixed reader to return appropriate string for formula cell
#
# Adds a new "standard format with formula inside" test case to both the
# "format" and "formatForceString" sheets: a header/label pair in column
# AC, plus a handful of formula cells (including one that references
# another formula cell in column AE, and one that is a quote-prefixed
# text value) used to exercise reading formula cells as strings.

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("format", "formatForceString")) {
    $ws = $wb.Worksheets.Item($sheetName)

    if ($sheetName -eq "format") {
        $ws.Range("AC1").Value = "format#standardAndFormulaInside"
    } else {
        $ws.Range("AC1").Value = "formatStr#standardAndFormulaInside?type=string"
    }

    $ws.Range("AC3").Value = "Format to be standard and formula inside"

    # Text formula -> "AAABBB"
    $ws.Range("AC4").Formula = '="AAA"&"BBB"'

    # Numeric formula -> 1
    $ws.Range("AC5").Formula = "=1"

    # Text formula (quoted literal) -> "1"
    $ws.Range("AC6").Formula = '="1"'

    # Formula referencing another formula-less numeric cell
    $ws.Range("AE7").Value = 1
    $ws.Range("AC7").Formula = "=AE7"

    # Formula referencing a quote-prefixed (text) "1" cell
    $ws.Range("AE8").Formula = "'1"
    $ws.Range("AC8").Formula = "=AE8"
}
